$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.958.42"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "3.067.68"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.17"
$ws.Range("E5").Value = "  -2.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.96"
$ws.Range("E6").Value = "  -2.48%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "3.056.07"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  +0.39%  "

$ws.Range("E10").Value = "  -3.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.14"
$ws.Range("E11").Value = "  -7.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.451"
$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000222"
$ws.Range("E13").Value = "  +3.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.12"
$ws.Range("E14").Value = "  -1.97%  "

$ws.Range("D15").Value = "3.561.22"
$ws.Range("E15").Value = "  -0.04%  "

$ws.Range("D16").Value = "62.935.55"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").Value = "3.069.09"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.61"
$ws.Range("E19").Value = "  +0.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.15"
$ws.Range("E20").Value = "  -4.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.29"
$ws.Range("E21").Value = "  -1.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.691"
$ws.Range("E22").Value = "  -0.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.08"
$ws.Range("E23").Value = "  -1.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.85"
$ws.Range("E24").Value = "  +2.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.05"
$ws.Range("E25").Value = "  -1.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.70"
$ws.Range("E27").Value = "  -1.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.07"
$ws.Range("E28").Value = "  -1.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.90"
$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("E31").Value = "  -6.61%  "

$ws.Range("E32").Value = "  +1.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.35"
$ws.Range("E33").Value = "  -5.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "56.81"
$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("E35").Value = "  +4.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.00"
$ws.Range("E36").Value = "  +2.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "477.52"
$ws.Range("E37").Value = "  -10.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0393"
$ws.Range("E38").Value = "  -3.82%  "

$ws.Range("D39").Value = "3.083.94"
$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0793"
$ws.Range("E40").Value = "  +0.69%  "

$ws.Range("E41").Value = "  -1.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.06"
$ws.Range("E42").Value = "  +0.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.60"
$ws.Range("E43").Value = "  +0.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.252"
$ws.Range("E44").Value = "  +0.26%  "

$ws.Range("D46").Value = "0.0₃0538"
$ws.Range("E46").Value = "  +11.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.36"
$ws.Range("E47").Value = "  +0.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.01"
$ws.Range("E48").Value = "  -1.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.34"
$ws.Range("E49").Value = "  +2.05%  "

$ws.Range("E50").Value = "  +2.18%  "

$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.28"
$ws.Range("E51").Value = "  +2.35%  "

